$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-03 Saturday" "2024-08-04 Sunday"

Replace-Text "461×8=" "790×7="
Replace-Text "403×9=" "910×2="
Replace-Text "669×9=" "192×9="
Replace-Text "870×7=" "250×3="
Replace-Text "606×5=" "754×7="
Replace-Text "976×6=" "216×8="
Replace-Text "463×4=" "332×6="
Replace-Text "826×9=" "404×5="
Replace-Text "730×7=" "248×6="
Replace-Text "322×3=" "678×4="
Replace-Text "961×9=" "851×3="
Replace-Text "685×2=" "601×5="
Replace-Text "759×7=" "713×5="
Replace-Text "647×2=" "382×4="
Replace-Text "311×6=" "956×9="
Replace-Text "925×7=" "196×6="
Replace-Text "545×4=" "702×2="
Replace-Text "922×7=" "887×8="
Replace-Text "617×9=" "524×2="
Replace-Text "527×4=" "453×8="
Replace-Text "212×7=" "316×3="
Replace-Text "434×6=" "616×6="
Replace-Text "488×7=" "652×3="
Replace-Text "257×2=" "606×9="
Replace-Text "946×7=" "654×6="
